# feat: add 2022-Q3 data
#
# Before: "总计" (totals), "2022-Q2", "2022-Q1"
# After:  "总计" (totals), "2022-Q3" (new), "2022-Q2", "2022-Q1"
#
# Inserts a new "2022-Q3" sheet (cloned from "2022-Q2" so it keeps identical
# formatting/layout) populated with the new quarter's fund data, and updates
# the "总计" summary sheet with a new row for 2022-Q3 (pushing the existing
# 2022-Q2 / 2022-Q1 rows down by one).

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q3" sheet -------------------------------------
# Copy "2022-Q2" (placing the copy right before it) so the new sheet inherits
# the exact same column layout / cell styles, then rename + restamp values.
$wb.Worksheets.Item("2022-Q2").Copy($wb.Worksheets.Item("2022-Q2"))
$wb.Worksheets.Item("2022-Q2 (2)").Name = "2022-Q3"

# Re-fetch sheet handles by name now that the sheet collection has changed
# (handles captured before an Add/Copy can end up pointing at the wrong
# worksheet once positions shift).
$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q1Sheet = $wb.Worksheets.Item("2022-Q1")

# New quarter's fund figures (everything else on the sheet is unchanged).
# These columns hold numeric-looking text (matches the source workbook, which
# stores them as strings, not numbers) - a leading apostrophe enters them as
# text, then resetting to the "Normal" style drops the resulting quote-prefix
# flag so the cell format matches the rest of the sheet.
$q3Sheet.Range("D2").Value = "'0.66"
$q3Sheet.Range("E2").Value = "'86.08"
$q3Sheet.Range("F2").Value = "'4.92"
$q3Sheet.Range("G2").Value = "'0.0325"
$q3Sheet.Range("D2:G2").Style = "Normal"
$q3Sheet.Range("H2").Value = 4

# --- 2. Update the "总计" (totals) sheet ------------------------------------
# Row 2 now refers to the new quarter, row 3 shifts to what used to be row 2's
# label (value stays 0.03), and a new row 4 is appended holding the old row
# 3's data (2022-Q1, value 0.02).
$totalSheet.Range("B2").Value = "2022-Q3"

$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("D3").Value = 0.03

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.02

# --- 3. Restore the original active sheet -----------------------------------
# "2022-Q1" was the selected tab before the edit; keep it that way.
$q1Sheet.Activate()
